$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 156 (the "誰が死ぬ？パズル" post), shifting all
# subsequent rows up by one.
$ws.Rows.Item(156).Delete()
